$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("A11").Value = "Energia e ambiente, Economia regionale e mercato immobiliare"
$ws.Range("B11").Value = "Q, R"
$ws.Range("C11").Value = 0.04681647940074907

# Row 12
$ws.Range("B12").Value = "E01, D03, P33"
$ws.Range("C12").Value = 0.01029962546816479

# Row 22
$ws.Range("A22").Value = "Energia e ambiente, Economia regionale e mercato immobiliare"
$ws.Range("B22").Value = "Q, R"
$ws.Range("C22").Value = 0.05540935672514619

# Row 23
$ws.Range("B23").Value = "A2, B30, A23, Z1, B26, A1, B20, E00, D92, B17, B1, E71, Z31, Z3, Z10, B32, B31, B00, B27, P5, P16, D91, Z13"
$ws.Range("C23").Value = 0.04346978557504873
